$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1366.5
$ws.Range("I70").Value = 1049.75
$ws.Range("K70").Value = 3149.25
$ws.Range("M70").Value = -2879.25
$ws.Range("H73").Value = 1366.5
$ws.Range("I73").Value = 1049.75
$ws.Range("K73").Value = 3149.25
$ws.Range("M73").Value = -2213.25
$ws.Range("H88").Value = 1441.25
$ws.Range("J88").Value = 1588.6666
$ws.Range("L88").Value = 1588.6666
$ws.Range("N88").Value = -2400.6666
$ws.Range("H91").Value = 1441.25
$ws.Range("J91").Value = 1588.6666
$ws.Range("L91").Value = 1588.6666
$ws.Range("N91").Value = -4396.6666
$ws.Range("H138").Value = 5160.722
$ws.Range("J138").Value = 5199.5884
$ws.Range("L138").Value = 15598.7652
$ws.Range("N138").Value = -25878.7652

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15422.538
$ws.Range("I32").Value = 15772.091
$ws.Range("K32").Value = 15772.091
$ws.Range("M32").Value = -15485.091
$ws.Range("H122").Value = 1667.3077
$ws.Range("I122").Value = 1582.5454
$ws.Range("K122").Value = 4747.6362
$ws.Range("M122").Value = -2297.6362

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8562
$ws.Range("I20").Value = 7083
$ws.Range("K20").Value = 7083
$ws.Range("M20").Value = -6836
$ws.Range("H86").Value = 3300.8823
$ws.Range("I86").Value = 3255.4614
$ws.Range("J86").Value = 3448.5
$ws.Range("K86").Value = 3255.4614
$ws.Range("L86").Value = 3448.5
$ws.Range("M86").Value = -2132.4614
$ws.Range("N86").Value = -5694.5
$ws.Range("H89").Value = 3300.8823
$ws.Range("I89").Value = 3255.4614
$ws.Range("J89").Value = 3448.5
$ws.Range("K89").Value = 16277.307
$ws.Range("L89").Value = 17242.5
$ws.Range("M89").Value = -10661.307
$ws.Range("N89").Value = -28474.5
$ws.Range("H134").Value = 4999
$ws.Range("I134").Value = 4999
$ws.Range("K134").Value = 14997
$ws.Range("M134").Value = -12462

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2393.1853
$ws.Range("I31").Value = 2425.48
$ws.Range("K31").Value = 2425.48
$ws.Range("M31").Value = -2130.48
$ws.Range("H34").Value = 2393.1853
$ws.Range("I34").Value = 2425.48
$ws.Range("K34").Value = 2425.48
$ws.Range("M34").Value = -2223.48
$ws.Range("H86").Value = 24541.723
$ws.Range("J86").Value = 38075.223
$ws.Range("L86").Value = 38075.223
$ws.Range("N86").Value = -40321.223
$ws.Range("H89").Value = 24541.723
$ws.Range("J89").Value = 38075.223
$ws.Range("L89").Value = 190376.115
$ws.Range("N89").Value = -201608.115
$ws.Range("H99").Value = 9562.143
$ws.Range("I99").Value = 9484
$ws.Range("K99").Value = 9484
$ws.Range("M99").Value = -7986
$ws.Range("H126").Value = 9562.143
$ws.Range("I126").Value = 9484
$ws.Range("K126").Value = 28452
$ws.Range("M126").Value = -25982
$ws.Range("H132").Value = 4996.3335
$ws.Range("I132").Value = 4995
$ws.Range("K132").Value = 14985
$ws.Range("M132").Value = -12455
$ws.Range("H134").Value = 3464.5715
$ws.Range("I134").Value = 3208.8333
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 9626.499899999999
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -7091.499899999999
$ws.Range("N134").Value = -20067
$ws.Range("H141").Value = 60445.047
$ws.Range("J141").Value = 58467.3
$ws.Range("L141").Value = 58467.3
$ws.Range("N141").Value = -68827.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2992.5
$ws.Range("J92").Value = 2992.3333
$ws.Range("L92").Value = 8976.999899999999
$ws.Range("N92").Value = -11472.9999
$ws.Range("H132").Value = 4887.2104
$ws.Range("I132").Value = 4849.5
$ws.Range("J132").Value = 4891.647
$ws.Range("K132").Value = 43645.5
$ws.Range("L132").Value = 44024.823
$ws.Range("M132").Value = -41115.5
$ws.Range("N132").Value = -49084.823

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15000000
$ws.Range("J11").Value = 5000000
$ws.Range("L11").Value = 5000000
$ws.Range("N11").Value = -5000278
$ws.Range("H43").Value = 14672.667
$ws.Range("J43").Value = 18259
$ws.Range("L43").Value = 18259
$ws.Range("N43").Value = -18561
$ws.Range("H80").Value = 7841.8335
$ws.Range("I80").Value = 7442.5
$ws.Range("J80").Value = 8041.5
$ws.Range("K80").Value = 7442.5
$ws.Range("L80").Value = 8041.5
$ws.Range("M80").Value = -6444.5
$ws.Range("N80").Value = -10037.5
$ws.Range("H83").Value = 7841.8335
$ws.Range("I83").Value = 7442.5
$ws.Range("J83").Value = 8041.5
$ws.Range("K83").Value = 37212.5
$ws.Range("L83").Value = 40207.5
$ws.Range("M83").Value = -32220.5
$ws.Range("N83").Value = -50191.5
$ws.Range("H122").Value = 1470.9
$ws.Range("I122").Value = 1470.9
$ws.Range("K122").Value = 4412.700000000001
$ws.Range("M122").Value = -1962.700000000001
$ws.Range("H126").Value = 5238.375
$ws.Range("J126").Value = 6465
$ws.Range("L126").Value = 19395
$ws.Range("N126").Value = -24335

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 431.66666
$ws.Range("I55").Value = 417.83334
$ws.Range("J55").Value = 459.33334
$ws.Range("K55").Value = 417.83334
$ws.Range("L55").Value = 459.33334
$ws.Range("M55").Value = -244.83334
$ws.Range("N55").Value = -805.33334
$ws.Range("H61").Value = 3485
$ws.Range("I61").Value = 3400.4167
$ws.Range("K61").Value = 3400.4167
$ws.Range("M61").Value = -3198.4167
$ws.Range("H68").Value = 2736.875
$ws.Range("J68").Value = 3332.3333
$ws.Range("L68").Value = 3332.3333
$ws.Range("N68").Value = -4830.3333
$ws.Range("H71").Value = 2736.875
$ws.Range("J71").Value = 3332.3333
$ws.Range("L71").Value = 16661.6665
$ws.Range("N71").Value = -24149.6665
$ws.Range("H113").Value = 3485
$ws.Range("I113").Value = 3400.4167
$ws.Range("K113").Value = 3400.4167
$ws.Range("M113").Value = -1230.4167
$ws.Range("H136").Value = 7898.875
$ws.Range("J136").Value = 8000
$ws.Range("L136").Value = 24000
$ws.Range("N136").Value = -29100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 64616.125
$ws.Range("I62").Value = 102266.2
$ws.Range("J62").Value = 1866
$ws.Range("K62").Value = 102266.2
$ws.Range("L62").Value = 1866
$ws.Range("M62").Value = -101642.2
$ws.Range("N62").Value = -3114
$ws.Range("H65").Value = 64616.125
$ws.Range("I65").Value = 102266.2
$ws.Range("J65").Value = 1866
$ws.Range("K65").Value = 511331
$ws.Range("L65").Value = 9330
$ws.Range("M65").Value = -508211
$ws.Range("N65").Value = -15570
$ws.Range("H113").Value = 1747.3
$ws.Range("I113").Value = 2945
$ws.Range("K113").Value = 8835
$ws.Range("M113").Value = -6665
$ws.Range("H126").Value = 2264.4211
$ws.Range("I126").Value = 2264.4211
$ws.Range("K126").Value = 6793.263300000001
$ws.Range("M126").Value = -4323.263300000001
